$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.698.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "'2.492.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'586.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'176.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.34%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.515"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.140"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.32%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "'0.338"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("D12").Value = "'4.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "'2.941.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "'25.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "'67.533.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "'0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "'2.499.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "'11.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "'7.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "'352.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "'4.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'70.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").Value = "'4.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").Value = "'1.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'9.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").Value = "'2.605.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "'0.0₃0915"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").Value = "'507.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "'7.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").Value = "'1.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "'1.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'0.122"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.72%  "
$ws.Range("D36").Value = "'161.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").Value = "'18.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").Value = "'18.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'1.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("D42").Value = "'0.330"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "'4.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("D44").Value = "'2.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("D45").Value = "'143.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.10%  "
$ws.Range("D46").Value = "'3.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").Value = "'0.0₆0261"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("D48").Value = "'0.514"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "'0.0744"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").Value = "'1.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").Value = "'0.588"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.23%  "
